# Generate Report for Archive
#
# Refresh the localization-status report: the files
#   5e6c4ed6-1ad7-4eb9-98f3-8be25b99713a.md
#   827905e9-2647-41d6-bef1-d5f85f0bdb4c.md
# have moved out of "Ready for handoff" and are now "In Translation".
# Update the Status column on every sheet that tracks it: the overview
# sheet (zh-cn/de-de status columns B & C) and each per-locale sheet
# (Status column C).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: rows 7 & 8 (5e6c4ed6..., 827905e9...), columns B (zh-cn) and C (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus
$overview.Range("B8").Value = $newStatus
$overview.Range("C8").Value = $newStatus

# --- Per-locale sheets: rows 7 & 8, column C (Status)
$locales = @("zh-cn", "de-de")
foreach ($locale in $locales) {
    $sheet = $wb.Worksheets.Item($locale)
    $sheet.Range("C7").Value = $newStatus
    $sheet.Range("C8").Value = $newStatus
}
